$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: each entry below is the new value for one
# changed cell (Coin/Link for the two rows that swapped position,
# and the refreshed Price / Volume(1h) text for every row whose
# market data moved). Column D (Price) and column E (Volume(1h))
# store plain text such as "29.958.05" or "  +0.59%  " as inline
# strings in the workbook, so each D/E cell is switched to the
# Text number format before its value is written -- otherwise
# Excel would parse strings like "1.002" or "241.80" as numbers
# and silently drop the thousands-style dot or the trailing zero.
$updates = @(
    @{ Row=2; D='29.958.05'; E='  +0.59%  ' }
    @{ Row=3; D='1.894.09'; E='  -0.26%  ' }
    @{ Row=4; D='1.002'; E='  +0.20%  ' }
    @{ Row=5; D='0.8214'; E='  +6.84%  ' }
    @{ Row=6; D='241.80'; E='  +0.67%  ' }
    @{ Row=7; D='1.003'; E='  +0.30%  ' }
    @{ Row=8; D='0.3228'; E='  +5.92%  ' }
    @{ Row=9; D='26.50'; E='  +4.46%  ' }
    @{ Row=10; D='0.07026'; E='  +2.83%  ' }
    @{ Row=11; D='0.08038'; E='  +0.76%  ' }
    @{ Row=12; D='0.7477'; E='  +1.52%  ' }
    @{ Row=13; D='1.905.18'; E='  +0.17%  ' }
    @{ Row=14; D='5.200'; E='  +0.78%  ' }
    @{ Row=15; D='92.34'; E='  +1.44%  ' }
    @{ Row=16; D='29.955.90'; E='  +0.53%  ' }
    @{ Row=17; D='14.04'; E='  +2.32%  ' }
    @{ Row=18; D='5.897'; E='  +0.32%  ' }
    @{ Row=19; D='244.86'; E='  +0.15%  ' }
    @{ Row=20; D='0.000007749'; E='  +0.71%  ' }
    @{ Row=21; D='1.002'; E='  +0.23%  ' }
    @{ Row=22; D='2.139.40'; E='  -1.16%  ' }
    @{ Row=23; D='1.002'; E='  +0.21%  ' }
    @{ Row=24; D='6.908'; E='  +0.45%  ' }
    @{ Row=25; D='0.1587'; E='  +23.51%  ' }
    @{ Row=26; D='166.69'; E='  -0.09%  ' }
    @{ Row=27; D='9.193'; E='  -0.49%  ' }
    @{ Row=28; E='  +1.04%  ' }
    @{ Row=29; D='2.069'; E='  +2.10%  ' }
    @{ Row=30; D='1.371'; E='  -1.98%  ' }
    @{ Row=31; D='1.517'; E='  +0.55%  ' }
    @{ Row=32; D='4.266'; E='  -0.06%  ' }
    @{ Row=33; D='0.05612'; E='  +7.02%  ' }
    @{ Row=34; D='4.076'; E='  +0.36%  ' }
    @{ Row=35; D='1.273'; E='  +2.66%  ' }
    @{ Row=36; D='0.7305' }
    @{ Row=37; D='2.725'; E='  +0.27%  ' }
    @{ Row=38; D='0.01914'; E='  +0.21%  ' }
    @{ Row=39; D='2.784'; E='  +0.22%  ' }
    @{ Row=40; D='0.4416'; E='  +0.32%  ' }
    @{ Row=41; D='71.95'; E='  +0.04%  ' }
    @{ Row=42; D='5.956'; E='  -3.76%  ' }
    @{ Row=43; D='0.8436'; E='  +1.18%  ' }
    @{ Row=44; D='1.003'; E='  +0.32%  ' }
    @{ Row=45; E='  +0.29%  ' }
    @{ Row=46; D='7.582'; E='  +0.10%  ' }
    @{ Row=47; D='100.68'; E='  +0.76%  ' }
    @{ Row=48; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='990.74'; E='  +9.29%  ' }
    @{ Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.676'; E='  -0.40%  ' }
    @{ Row=50; D='2.043.16'; E='  -0.59%  ' }
    @{ Row=51; D='36.01'; E='  -0.40%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B$($u.Row)").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$($u.Row)").Value = $u.C }
    if ($u.ContainsKey("D")) {
        $ws.Range("D$($u.Row)").NumberFormat = "@"
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$($u.Row)").NumberFormat = "@"
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
